$d = $word.ActiveDocument

# The wml.xsd schema for CT_RPr requires <w:rPr> child elements to appear
# in a fixed order (rStyle, rFonts, b, bCs, i, iCs, ... color, ...).
# Several custom character styles in styles.xml had <w:color> emitted
# before <w:b>/<w:i>, which OOXMLValidatorCLI flags as a schema error
# even though xmllint stays silent. Re-assert the existing bold/italic
# flags on each affected style so the run-properties get re-serialized
# in the schema-correct order (b/i before color).

$boldStyles = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)
foreach ($name in $boldStyles) {
    $s = $d.Styles($name)
    $s.Font.Bold = $s.Font.Bold
}

$italicStyles = @(
    "CommentTok",
    "DocumentationTok"
)
foreach ($name in $italicStyles) {
    $s = $d.Styles($name)
    $s.Font.Italic = $s.Font.Italic
}

$boldItalicStyles = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)
foreach ($name in $boldItalicStyles) {
    $s = $d.Styles($name)
    $s.Font.Bold = $s.Font.Bold
    $s.Font.Italic = $s.Font.Italic
}
